# Add a new date column (AE) for 16-10-2020 with that day's cumulative
# case counts for every State/UT, mirroring the formatting used by the
# preceding date-header column (AD1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (row 1) ---------------------------------------------
$headerCell = $ws.Range("AE1")
$headerCell.Value = "16-10-2020"

# Match the look of the neighbouring date headers (N1:AD1): bold text,
# centered/top aligned, thin box border all around.
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108  # xlCenter
$headerCell.VerticalAlignment = -4160    # xlTop
$headerCell.Borders.LineStyle = 1        # xlContinuous
$headerCell.Borders.Weight = 2           # xlThin

# --- New daily totals for each State/UT (row -> value) ---------------
$values = @{
    2  = 3817
    3  = 725099
    4  = 9889
    5  = 169335
    6  = 189186
    7  = 12232
    8  = 123943
    9  = 3099
    10 = 292502
    11 = 35161
    12 = 137733
    13 = 134719
    14 = 15389
    15 = 75641
    16 = 87240
    17 = 620008
    18 = 222231
    19 = 4310
    20 = 139717
    21 = 1330483
    22 = 11081
    23 = 5646
    24 = 2121
    25 = 6017
    26 = 238535
    27 = 27365
    28 = 115186
    29 = 143984
    30 = 3129
    31 = 622458
    32 = 194653
    33 = 25765
    34 = 49997
    35 = 404545
    36 = 271563
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 31).Value = $values[$row]
}
